# Translate English month names to Portuguese on the "por mês" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("por mês")

$months = @("Janeiro", "Fevereiro", "Março", "Abril", "Maio", "Junho", "Julho", "Agosto", "Setembro", "Outubro", "Novembro", "Dezembro")

for ($i = 0; $i -lt $months.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $months[$i]
}

$wb.Save()
